$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("outputs/2024-04-19/06-09-19", $false, "eicu",     "mlm",         "", "descemb_bert", "NV", $null),
    @("outputs/2024-04-19/06-10-24", $false, "mimiciii", "mlm",         "", "descemb_bert", "NV", $null),
    @("outputs/2024-04-19/06-11-39", $false, "mimiciii", "mlm",         "", "descemb_bert", "NV", 0),
    @("outputs/2024-04-19/06-16-24", $false, "mimiciii", "readmission", "", "descemb_bert", "NV", $null),
    @("outputs/2024-04-19/06-16-52", $false, "mimiciii", "mlm",         "", "descemb_bert", "NV", $null),
    @("outputs/2024-04-19/06-18-04", $true,  "mimiciii", "mlm",         "", "descemb_bert", "NV", 0),
    @("outputs/2024-04-19/15-55-40", $true,  "mimiciii", "mlm",         "", "descemb_rnn",  "NV", 0)
)

$startRow = 139
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    if ($row[7] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[7]
    }
}
